$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format before writing, so numeric-looking
# strings (e.g. "593.38", "0.173") are stored as text, matching the
# original inline-string cells instead of being coerced to numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.509.76"
$ws.Range("E2").Value = "  -4.51%  "
$ws.Range("D3").Value = "3.284.08"
$ws.Range("E3").Value = "  -7.04%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "593.38"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "152.72"
$ws.Range("E6").Value = "  -11.41%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.273.20"
$ws.Range("E8").Value = "  -7.29%  "
$ws.Range("E9").Value = "  -10.44%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -13.30%  "
$ws.Range("D11").Value = "6.74"
$ws.Range("E11").Value = "  -6.82%  "
$ws.Range("D12").Value = "0.514"
$ws.Range("E12").Value = "  -12.42%  "
$ws.Range("D13").Value = "38.86"
$ws.Range("E13").Value = "  -15.95%  "
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  -10.77%  "
$ws.Range("D15").Value = "3.815.01"
$ws.Range("E15").Value = "  -7.01%  "
$ws.Range("D16").Value = "67.559.62"
$ws.Range("E16").Value = "  -4.66%  "
$ws.Range("D17").Value = "3.288.40"
$ws.Range("E17").Value = "  -6.91%  "
$ws.Range("D18").Value = "7.29"
$ws.Range("E18").Value = "  -13.56%  "
$ws.Range("D19").Value = "536.68"
$ws.Range("E19").Value = "  -11.63%  "
$ws.Range("E20").Value = "  -6.19%  "
$ws.Range("D21").Value = "15.18"
$ws.Range("E21").Value = "  -14.19%  "
$ws.Range("D22").Value = "0.765"
$ws.Range("E22").Value = "  -13.20%  "
$ws.Range("D23").Value = "7.91"
$ws.Range("E23").Value = "  -13.01%  "
$ws.Range("D24").Value = "13.72"
$ws.Range("E24").Value = "  -12.04%  "
$ws.Range("D25").Value = "86.03"
$ws.Range("E25").Value = "  -11.84%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "3.27"
$ws.Range("E27").Value = "  -11.69%  "
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  -10.02%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "29.50"
$ws.Range("E29").Value = "  -12.35%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -15.73%  "
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  -10.11%  "
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -10.87%  "
$ws.Range("D33").Value = "547.69"
$ws.Range("E33").Value = "  -11.05%  "
$ws.Range("D34").Value = "6.66"
$ws.Range("E34").Value = "  -17.97%  "
$ws.Range("D35").Value = "5.82"
$ws.Range("E35").Value = "  -14.49%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.0470"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "53.57"
$ws.Range("E38").Value = "  -5.68%  "
$ws.Range("D39").Value = "0.0867"
$ws.Range("E39").Value = "  -12.77%  "
$ws.Range("D40").Value = "9.17"
$ws.Range("E40").Value = "  -15.43%  "
$ws.Range("E41").Value = "  -9.41%  "
$ws.Range("E42").Value = "  -17.76%  "
$ws.Range("D43").Value = "2.950.74"
$ws.Range("E43").Value = "  -11.88%  "
$ws.Range("D44").Value = "0.271"
$ws.Range("E44").Value = "  -12.59%  "
$ws.Range("D45").Value = "0.0₃0600"
$ws.Range("E45").Value = "  -17.50%  "
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").Value = "  -11.01%  "
$ws.Range("E47").Value = "  -15.02%  "
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  -17.90%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "126.63"
$ws.Range("E50").Value = "  -5.55%  "
$ws.Range("E51").Value = "  -12.05%  "

# Restore default (General/Normal) styling so the format override above
# doesn't leave a lingering style on these cells.
$priceRange.Style = "Normal"
